# Updated symbol list with GitHub Actions - refresh Price (column D) values
# and a few Volume(1h)/label (column E) cells to the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$Text)
    $range = $Worksheet.Range($Address)
    $originalStyle = $range.Style
    # Force the cell into text mode while writing so Excel doesn't coerce a
    # numeric-looking string (e.g. "244.50") into a number and drop
    # significant trailing zeros, then restore the original cell style.
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = $originalStyle
}

# Column D ("Price") updates
Set-TextValue $ws "D2"  "244.50"
Set-TextValue $ws "D3"  "23.52"
Set-TextValue $ws "D4"  "5.729"
Set-TextValue $ws "D5"  "0.05830"
Set-TextValue $ws "D6"  "3.416"
Set-TextValue $ws "D7"  "6.466"
Set-TextValue $ws "D9"  "0.7985"
Set-TextValue $ws "D10" "0.1466"
Set-TextValue $ws "D11" "0.07617"
Set-TextValue $ws "D12" "0.03225"
Set-TextValue $ws "D13" "0.02963"
Set-TextValue $ws "D14" "0.09232"
Set-TextValue $ws "D15" "0.001657"
Set-TextValue $ws "D16" "3.265"
Set-TextValue $ws "D17" "0.04753"
Set-TextValue $ws "D18" "0.0005981"
Set-TextValue $ws "D19" "0.006271"
Set-TextValue $ws "D20" "0.005415"
Set-TextValue $ws "D21" "0.001063"
Set-TextValue $ws "D22" "0.0001498"
Set-TextValue $ws "D23" "3.696"
Set-TextValue $ws "D24" "2.192"
Set-TextValue $ws "D25" "0.3342"
Set-TextValue $ws "D27" "0.0009981"
Set-TextValue $ws "D40" "0.04297"
Set-TextValue $ws "D41" "0.007073"
Set-TextValue $ws "D42" "0.003594"
Set-TextValue $ws "D43" "0.1062"
Set-TextValue $ws "D44" "0.009748"
Set-TextValue $ws "D46" "0.00005428"
Set-TextValue $ws "D47" "0.00000000749"
Set-TextValue $ws "D48" "0.7840"
Set-TextValue $ws "D49" "0.1018"
Set-TextValue $ws "D50" "0.00002096"
Set-TextValue $ws "D51" "0.01008"

# Column E ("Volume(1h)") label updates
Set-TextValue $ws "E16" "15MCDexMCB"
Set-TextValue $ws "E18" "17OneONEWorstin24h"
Set-TextValue $ws "E49" "48BOLOBOLOBestin24h"
